# Applies the "Updated cryptos list" data refresh described by the commit diff.
# Source rows are keyed by worksheet row number (row 2 = first coin, etc.).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price values (column D) are plain numeric-looking strings (e.g. "0.529", "613.17").
# Force those cells to Text format first so Excel keeps them as literal strings
# (matching the original inlineStr formatting) instead of silently converting them
# to numbers (which would drop trailing zeros / change precision).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# --- Cell value updates ---
$ws.Range("D2").Value = "69.628.09"
$ws.Range("E2").Value = "  -0.68%  "
$ws.Range("D3").Value = "3.795.05"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "613.17"
$ws.Range("E5").Value = "  -1.70%  "
$ws.Range("D6").Value = "177.22"
$ws.Range("E6").Value = "  -1.61%  "
$ws.Range("D7").Value = "3.792.39"
$ws.Range("E7").Value = "  +0.78%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "0.529"
$ws.Range("E9").Value = "  -1.20%  "
$ws.Range("E10").Value = "  -1.38%  "
$ws.Range("D11").Value = "6.48"
$ws.Range("E11").Value = "  +2.41%  "
$ws.Range("E12").Value = "  -1.37%  "
$ws.Range("D13").Value = "39.77"
$ws.Range("E13").Value = "  -3.58%  "
$ws.Range("E14").Value = "  -2.51%  "
$ws.Range("D15").Value = "4.431.61"
$ws.Range("E15").Value = "  +1.08%  "
$ws.Range("D16").Value = "3.800.32"
$ws.Range("E16").Value = "  +1.11%  "
$ws.Range("D17").Value = "69.695.42"
$ws.Range("E17").Value = "  -0.68%  "
$ws.Range("D18").Value = "7.55"
$ws.Range("E18").Value = "  -1.57%  "
$ws.Range("E19").Value = "  -4.01%  "
$ws.Range("D20").Value = "16.62"
$ws.Range("E20").Value = "  -1.04%  "
$ws.Range("D21").Value = "506.93"
$ws.Range("E21").Value = "  -0.22%  "
$ws.Range("E22").Value = "  +1.66%  "
$ws.Range("D23").Value = "0.735"
$ws.Range("E23").Value = "  +0.67%  "
$ws.Range("D24").Value = "2.48"
$ws.Range("E24").Value = "  -1.53%  "
$ws.Range("D25").Value = "86.27"
$ws.Range("E25").Value = "  -1.17%  "
$ws.Range("E26").Value = "  +4.79%  "
$ws.Range("D27").Value = "12.63"
$ws.Range("E27").Value = "  -4.61%  "
$ws.Range("D28").Value = "10.56"
$ws.Range("E28").Value = "  -5.36%  "
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("D30").Value = "2.52"
$ws.Range("E30").Value = "  +0.53%  "
$ws.Range("E31").Value = "  +0.55%  "
$ws.Range("D32").Value = "8.04"
$ws.Range("E32").Value = "  +1.56%  "
$ws.Range("D33").Value = "31.41"
$ws.Range("E33").Value = "  -0.54%  "
$ws.Range("E34").Value = "  -1.82%  "
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").Value = "1.04"
$ws.Range("E36").Value = "  -2.51%  "
$ws.Range("D37").Value = "6.11"
$ws.Range("E37").Value = "  -2.21%  "
$ws.Range("E38").Value = "  +6.16%  "
$ws.Range("D39").Value = "480.94"
$ws.Range("E39").Value = "  +13.31%  "
$ws.Range("D40").Value = "0.338"
$ws.Range("E40").Value = "  +0.56%  "
$ws.Range("D41").Value = "3.04"
$ws.Range("E41").Value = "  +6.02%  "
$ws.Range("D42").Value = "2.06"
$ws.Range("E42").Value = "  -2.85%  "
$ws.Range("E43").Value = "  -1.22%  "
$ws.Range("D44").Value = "44.16"
$ws.Range("E44").Value = "  -2.37%  "
$ws.Range("D45").Value = "8.57"
$ws.Range("E45").Value = "  -2.24%  "
$ws.Range("D46").Value = "2.930.63"
$ws.Range("E46").Value = "  -2.64%  "
$ws.Range("D47").Value = "0.0362"
$ws.Range("E47").Value = "  -0.85%  "
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").Value = "139.24"
$ws.Range("E48").Value = "  +0.36%  "
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "27.24"
$ws.Range("E50").Value = "  -0.80%  "
$ws.Range("E51").Value = "  -4.19%  "

# Restore default styling on the cells we forced to Text format, so only the
# values change and no stray cell formatting is introduced.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"
